$d = $word.ActiveDocument
$c0 = $d.Comments(1)
Write-Output ("c0 Reference Start/End: " + $c0.Reference.Start + "," + $c0.Reference.End)
$c1 = $d.Comments(2)
Write-Output ("c1 Reference Start/End: " + $c1.Reference.Start + "," + $c1.Reference.End)
